$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 content updates -------------------------------------------------
# TaskEstimatedTime: 5 -> 2
$ws.Range("N2").Value = 2
# TaskStartDate: 2024-02-24 -> 2024-02-27
$ws.Range("O2").Value = "2024-02-27"

# Normalize Q2's formatting (drop the redundant number-format override so it
# matches the plain bordered style used elsewhere in the row).
$ws.Range("A2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Hyperlinks --------------------------------------------------------
# Drop every hyperlink on the sheet (there are two: D2 and D3) and re-create
# only the one that should survive (D2), since row 3 (and its hyperlink) is
# being removed below.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:john@gmail.com")
# Re-touch the border so the cell resolves back to the shared "Hyperlink"
# cell-style xf instead of a freshly minted duplicate.
$ws.Range("D2").Borders.LineStyle = 1

# --- Remove the second (defect) test case row ------------------------------
$ws.Rows("3").Delete()

# --- View/selection state ---------------------------------------------
$ws.Range("A2").Select()
